$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E6 end time (10:50pm instead of 10:50am)
$ws.Range("E6").Value = 0.95138888888888884

# Fill in row 7 with a new log entry
$ws.Range("B7").Value = 9417
$ws.Range("C7").Value = "30/03/2020"
$ws.Range("D7").Value = "8:30pm"
$ws.Range("E7").Value = 0.95833333333333337
$ws.Range("G7").Value = "Finished Logic Unit Design"

# Update the active selection
$ws.Range("E10").Select()
